{"js": "// The diff does two things to JustRun/README.docx:\n//  1. Removes the stray \"_GoBack\" bookmark that sat at the end of the\n//     \"Features:\" paragraph.\n//  2. Appends a new chunk of commentary text to the end of the\n//     \"In Past Workouts, ...\" paragraph, with the \"_GoBack\" bookmark\n//     re-inserted in the middle of that new text (right after\n//     \"...issues w\", before \"ith TabBarController...\").\n\n// --- 1. Drop the old _GoBack bookmark -------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2. Append the new commentary sentence --------------------------------\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text && p.text.indexOf(\"In Past Workouts\") === 0\n);\nif (!target) {\n  throw new Error(\"Could not find the 'In Past Workouts' paragraph\");\n}\n\n// First chunk (up through \"...issues w\") and the remainder, split so the\n// bookmark can be re-inserted exactly between them, mirroring the XML diff\n// (<w:bookmarkStart/> / <w:bookmarkEnd/> sit between the two runs).\nconst firstChunk = \" **We ran into issues w\";\nconst restChunk =\n  \"ith TabBarController reinitializing views every time a tab was clicked, \" +\n  \"which caused problems with delegating data between the RunViewController \" +\n  \"and PastWorkoutsController. \";\n\n// Insert all of the new text first (appended at the end of the paragraph).\ntarget.insertText(firstChunk + restChunk, \"End\");\nawait context.sync();\n\n// Now find the boundary between the two chunks inside the paragraph and drop\n// the bookmark there, using a fresh search hit (not a stale cached Range) so\n// the insertion point is resolved against the paragraph's current content.\nconst hits = target.search(firstChunk, { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length) {\n  const boundary = hits.items[hits.items.length - 1].getRange(\"After\");\n  boundary.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The diff does two things to JustRun/README.docx:\n#  1. Removes the stray \"_GoBack\" bookmark that sat at the end of the\n#     \"Features:\" paragraph.\n#  2. Appends a new chunk of commentary text to the end of the\n#     \"In Past Workouts, ...\" paragraph, with the \"_GoBack\" bookmark\n#     re-inserted in the middle of that new text (right after\n#     \"...issues w\", before \"ith TabBarController...\").\n\n$d = $word.ActiveDocument\n\n# --- 1. Drop the old _GoBack bookmark ---------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- 2. Append the new commentary sentence -----------------------------\n$firstChunk = \" **We ran into issues w\"\n$restChunk = \"ith TabBarController reinitializing views every time a tab was clicked, which caused problems with delegating data between the RunViewController and PastWorkoutsController. \"\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"In Past Workouts\")) {\n        $rng = $p.Range\n        $rng.Collapse(0)     # wdCollapseEnd\n        $rng.MoveEnd(1, -1)  # wdCharacter: back off the paragraph mark\n        $rng.InsertAfter($firstChunk + $restChunk)\n        break\n    }\n}\n\n# Re-insert the bookmark exactly between the two chunks: find the unique\n# boundary text and collapse to its end, then drop the bookmark there (using\n# a freshly-resolved Find range rather than a stale cached Range, so the\n# insertion point reflects the document's current state).\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$found = $findRange.Find.Execute($firstChunk)\nif ($found) {\n    $findRange.Collapse(0)  # wdCollapseEnd\n    $d.Bookmarks.Add(\"_GoBack\", $findRange)\n}\n"}
